$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 742.1579
$ws.Range("I33").Value = 163.72728
$ws.Range("J33").Value = 1537.5
$ws.Range("K33").Value = 163.72728
$ws.Range("L33").Value = 1537.5
$ws.Range("M33").Value = 65.27271999999999
$ws.Range("N33").Value = -1995.5
$ws.Range("H112").Value = 869372.25
$ws.Range("J112").Value = 1043113.4
$ws.Range("L112").Value = 3129340.2
$ws.Range("N112").Value = -3131556.2
$ws.Range("H113").Value = 2004.1666
$ws.Range("I113").Value = 1941.4706
$ws.Range("J113").Value = 2156.4285
$ws.Range("K113").Value = 1941.4706
$ws.Range("L113").Value = 2156.4285
$ws.Range("M113").Value = 1312.5294
$ws.Range("N113").Value = -8664.4285
$ws.Range("H141").Value = 114860.375
$ws.Range("I141").Value = 3998.75
$ws.Range("J141").Value = 225722
$ws.Range("K141").Value = 11996.25
$ws.Range("L141").Value = 677166
$ws.Range("M141").Value = -6816.25
$ws.Range("N141").Value = -687526
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2118.8853
$ws.Range("I2").Value = 448.275
$ws.Range("J2").Value = 5301
$ws.Range("K2").Value = 448.275
$ws.Range("L2").Value = 5301
$ws.Range("M2").Value = -335.275
$ws.Range("N2").Value = -5527
$ws.Range("H45").Value = 1226.9
$ws.Range("I45").Value = 1147.25
$ws.Range("J45").Value = 1280
$ws.Range("K45").Value = 1147.25
$ws.Range("L45").Value = 1280
$ws.Range("M45").Value = -770.25
$ws.Range("N45").Value = -2034
$ws.Range("H116").Value = 2118.8853
$ws.Range("I116").Value = 448.275
$ws.Range("J116").Value = 5301
$ws.Range("K116").Value = 448.275
$ws.Range("L116").Value = 5301
$ws.Range("M116").Value = 1845.725
$ws.Range("N116").Value = -9889
$ws.Range("H122").Value = 1276.9375
$ws.Range("I122").Value = 1225.2858
$ws.Range("K122").Value = 3675.8574
$ws.Range("M122").Value = -1225.8574
$ws.Range("H132").Value = 156916.77
$ws.Range("I132").Value = 22620.145
$ws.Range("J132").Value = 591405.8
$ws.Range("K132").Value = 67860.435
$ws.Range("L132").Value = 1774217.4
$ws.Range("M132").Value = -65330.435
$ws.Range("N132").Value = -1779277.4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2118.8853
$ws.Range("I3").Value = 448.275
$ws.Range("J3").Value = 5301
$ws.Range("K3").Value = 448.275
$ws.Range("L3").Value = 5301
$ws.Range("M3").Value = -334.275
$ws.Range("N3").Value = -5529
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 468.2143
$ws.Range("I5").Value = 49.333332
$ws.Range("J5").Value = 782.375
$ws.Range("K5").Value = 49.333332
$ws.Range("L5").Value = 782.375
$ws.Range("M5").Value = 62.666668
$ws.Range("N5").Value = -1006.375
$ws.Range("H122").Value = 3337.3333
$ws.Range("I122").Value = 3337.3333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 10011.9999
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = -7561.999899999999
$ws.Range("H132").Value = 64239
$ws.Range("I132").Value = 78148.16
$ws.Range("K132").Value = 234444.48
$ws.Range("M132").Value = -231914.48
$ws.Range("H134").Value = 1309.4651
$ws.Range("I134").Value = 1056.9429
$ws.Range("J134").Value = 2414.25
$ws.Range("K134").Value = 3170.8287
$ws.Range("L134").Value = 7242.75
$ws.Range("M134").Value = -635.8287
$ws.Range("N134").Value = -12312.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1438.3889
$ws.Range("I5").Value = 1962.3636
$ws.Range("J5").Value = 615
$ws.Range("K5").Value = 5887.0908
$ws.Range("L5").Value = 1845
$ws.Range("M5").Value = -5775.0908
$ws.Range("N5").Value = -2069
$ws.Range("H135").Value = 1438.3889
$ws.Range("I135").Value = 1962.3636
$ws.Range("J135").Value = 615
$ws.Range("K135").Value = 17661.2724
$ws.Range("L135").Value = 5535
$ws.Range("M135").Value = -15126.2724
$ws.Range("N135").Value = -10605
$ws.Range("H136").Value = 1658.5454
$ws.Range("I136").Value = 1637.625
$ws.Range("J136").Value = 1714.3334
$ws.Range("K136").Value = 4912.875
$ws.Range("L136").Value = 5143.0002
$ws.Range("M136").Value = 187.125
$ws.Range("N136").Value = -15343.0002
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 27308.334
$ws.Range("J42").Value = 27308.334
$ws.Range("L42").Value = 27308.334
$ws.Range("N42").Value = -28278.334
$ws.Range("H70").Value = 4345.8823
$ws.Range("I70").Value = 3884.2856
$ws.Range("J70").Value = 4669
$ws.Range("K70").Value = 3884.2856
$ws.Range("L70").Value = 4669
$ws.Range("M70").Value = -3614.2856
$ws.Range("N70").Value = -5209
$ws.Range("H73").Value = 4345.8823
$ws.Range("I73").Value = 3884.2856
$ws.Range("J73").Value = 4669
$ws.Range("K73").Value = 3884.2856
$ws.Range("L73").Value = 4669
$ws.Range("M73").Value = -2948.2856
$ws.Range("N73").Value = -6541
$ws.Range("H95").Value = 20000
$ws.Range("J95").Value = 20000
$ws.Range("L95").Value = 20000
$ws.Range("N95").Value = -25492
$ws.Range("H102").Value = 44695
$ws.Range("I102").Value = 23856
$ws.Range("K102").Value = 23856
$ws.Range("M102").Value = -22234
$ws.Range("H115").Value = 27308.334
$ws.Range("J115").Value = 27308.334
$ws.Range("L115").Value = 27308.334
$ws.Range("N115").Value = -29658.334
$ws.Range("H122").Value = 2267.8823
$ws.Range("I122").Value = 2968.5
$ws.Range("J122").Value = 1645.1111
$ws.Range("K122").Value = 8905.5
$ws.Range("L122").Value = 4935.3333
$ws.Range("M122").Value = -6455.5
$ws.Range("N122").Value = -9835.3333
$ws.Range("H130").Value = 25756
$ws.Range("J130").Value = 25756
$ws.Range("L130").Value = 25756
$ws.Range("N130").Value = -35796
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 61461.766
$ws.Range("I40").Value = 750
$ws.Range("J40").Value = 69556.664
$ws.Range("K40").Value = 750
$ws.Range("L40").Value = 69556.664
$ws.Range("M40").Value = -614
$ws.Range("N40").Value = -69828.664
$ws.Range("H46").Value = 1576.8518
$ws.Range("I46").Value = 745.25
$ws.Range("J46").Value = 1721.4783
$ws.Range("K46").Value = 745.25
$ws.Range("L46").Value = 1721.4783
$ws.Range("M46").Value = -557.25
$ws.Range("N46").Value = -2097.4783
$ws.Range("H122").Value = 2611.724
$ws.Range("I122").Value = 2673.913
$ws.Range("J122").Value = 2373.3333
$ws.Range("K122").Value = 8021.739
$ws.Range("L122").Value = 7119.999899999999
$ws.Range("M122").Value = -5571.739
$ws.Range("N122").Value = -12019.9999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4484.96
$ws.Range("I62").Value = 3373
$ws.Range("J62").Value = 4696.7617
$ws.Range("K62").Value = 3373
$ws.Range("L62").Value = 4696.7617
$ws.Range("M62").Value = -2749
$ws.Range("N62").Value = -5944.7617
$ws.Range("H65").Value = 4484.96
$ws.Range("I65").Value = 3373
$ws.Range("J65").Value = 4696.7617
$ws.Range("K65").Value = 16865
$ws.Range("L65").Value = 23483.8085
$ws.Range("M65").Value = -13745
$ws.Range("N65").Value = -29723.8085
$ws.Range("H107").Value = 269.9
$ws.Range("I107").Value = 285.18182
$ws.Range("J107").Value = 251.22223
$ws.Range("K107").Value = 855.54546
$ws.Range("L107").Value = 753.66669
$ws.Range("M107").Value = 1064.45454
$ws.Range("N107").Value = -4593.66669
$ws.Range("H123").Value = 25004.158
$ws.Range("J123").Value = 25004.158
$ws.Range("L123").Value = 25004.158
$ws.Range("N123").Value = -34804.158
$ws.Range("H125").Value = 28599.8
$ws.Range("J125").Value = 28599.8
$ws.Range("L125").Value = 28599.8
$ws.Range("N125").Value = -38439.8
$ws.Range("H136").Value = 1358884.8
$ws.Range("I136").Value = 1744127.6
$ws.Range("J136").Value = 527571.2
$ws.Range("K136").Value = 5232382.800000001
$ws.Range("L136").Value = 1582713.6
$ws.Range("M136").Value = -5229832.800000001
$ws.Range("N136").Value = -1587813.6
